$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I5").Value = 52.333332
$ws.Range("M5").Value = 62.666668
$ws.Range("K5").Value = 52.333332
$ws.Range("H5").Value = 68.111115
$ws.Range("J19").Value = 2000
$ws.Range("I19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 2000
$ws.Range("H19").Value = 2000
$ws.Range("N19").Value = -2350
$ws.Range("J51").Value = 0
$ws.Range("I51").Value = 6498
$ws.Range("M51").Value = -6014
$ws.Range("K51").Value = 6498
$ws.Range("L51").Value = 0
$ws.Range("H51").Value = 6498
$ws.Range("N51").ClearContents()
$ws.Range("I137").Value = 1282.7894
$ws.Range("M137").Value = -1298.3682
$ws.Range("K137").Value = 3848.3682
$ws.Range("H137").Value = 1684.1482
$ws.Range("I138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("K138").Value = 0
$ws.Range("H138").Value = 2000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 11694.3
$ws.Range("M32").Value = -11407.3
$ws.Range("K32").Value = 11694.3
$ws.Range("H32").Value = 13358.454
$ws.Range("I61").Value = 4420
$ws.Range("M61").Value = -4208
$ws.Range("K61").Value = 4420
$ws.Range("H61").Value = 4420
$ws.Range("I102").Value = 3110
$ws.Range("M102").Value = -1488
$ws.Range("K102").Value = 3110
$ws.Range("H102").Value = 3180
$ws.Range("I136").Value = 4420
$ws.Range("M136").Value = -10710
$ws.Range("K136").Value = 13260
$ws.Range("H136").Value = 4420

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I22").Value = 132.66667
$ws.Range("M22").Value = 40.33332999999999
$ws.Range("K22").Value = 132.66667
$ws.Range("H22").Value = 132.66667
$ws.Range("J54").Value = 42894.5
$ws.Range("I54").Value = 8586.888999999999
$ws.Range("M54").Value = -8102.888999999999
$ws.Range("K54").Value = 8586.888999999999
$ws.Range("L54").Value = 42894.5
$ws.Range("H54").Value = 14824.637
$ws.Range("N54").Value = -43862.5
$ws.Range("I99").Value = 1004.75
$ws.Range("M99").Value = 493.25
$ws.Range("K99").Value = 1004.75
$ws.Range("H99").Value = 1036.3334
$ws.Range("J135").Value = 39998.668
$ws.Range("L135").Value = 39998.668
$ws.Range("H135").Value = 39998.668
$ws.Range("N135").Value = -50138.668
$ws.Range("J137").Value = 78000
$ws.Range("L137").Value = 78000
$ws.Range("H137").Value = 78000
$ws.Range("N137").Value = -88200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I58").Value = 3606.375
$ws.Range("M58").Value = -3403.375
$ws.Range("K58").Value = 3606.375
$ws.Range("H58").Value = 5963.263
$ws.Range("J62").Value = 7000
$ws.Range("L62").Value = 7000
$ws.Range("H62").Value = 7000
$ws.Range("N62").Value = -8248
$ws.Range("J65").Value = 7000
$ws.Range("L65").Value = 35000
$ws.Range("H65").Value = 7000
$ws.Range("N65").Value = -41240
$ws.Range("I107").Value = 936.5454999999999
$ws.Range("M107").Value = 983.4545000000001
$ws.Range("K107").Value = 936.5454999999999
$ws.Range("H107").Value = 858.5714
$ws.Range("J132").Value = 2749.5
$ws.Range("L132").Value = 8248.5
$ws.Range("H132").Value = 2749.5
$ws.Range("N132").Value = -13308.5
$ws.Range("I134").Value = 2903.9167
$ws.Range("M134").Value = -6176.750100000001
$ws.Range("K134").Value = 8711.750100000001
$ws.Range("H134").Value = 3603.75
$ws.Range("I136").Value = 3606.375
$ws.Range("M136").Value = -8269.125
$ws.Range("K136").Value = 10819.125
$ws.Range("H136").Value = 5963.263

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J12").Value = 268.1111
$ws.Range("I12").Value = 90
$ws.Range("M12").Value = -97
$ws.Range("K12").Value = 270
$ws.Range("L12").Value = 804.3333
$ws.Range("H12").Value = 250.3
$ws.Range("N12").Value = -1150.3333
$ws.Range("J127").Value = 20000
$ws.Range("L127").Value = 60000
$ws.Range("H127").Value = 20000
$ws.Range("N127").Value = -69920
$ws.Range("I128").Value = 339894
$ws.Range("M128").Value = -1014702
$ws.Range("K128").Value = 1019682
$ws.Range("H128").Value = 339894
$ws.Range("I138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("K138").Value = 0
$ws.Range("H138").Value = 7307.769
$ws.Range("I139").Value = 3997.75
$ws.Range("M139").Value = -6853.25
$ws.Range("K139").Value = 11993.25
$ws.Range("H139").Value = 3997.75
$ws.Range("I140").Value = 1368.8182
$ws.Range("M140").Value = 1073.5454
$ws.Range("K140").Value = 4106.4546
$ws.Range("H140").Value = 1368.8182

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I2").Value = 211.07692
$ws.Range("M2").Value = -98.07692
$ws.Range("K2").Value = 211.07692
$ws.Range("H2").Value = 213.25
$ws.Range("I102").Value = 2000
$ws.Range("M102").Value = -378
$ws.Range("K102").Value = 2000
$ws.Range("H102").Value = 2000
$ws.Range("I126").Value = 1394.2
$ws.Range("M126").Value = -1712.6
$ws.Range("K126").Value = 4182.6
$ws.Range("H126").Value = 1661.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J46").Value = 4000
$ws.Range("I46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 4000
$ws.Range("H46").Value = 4000
$ws.Range("N46").Value = -4376
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("H56").Value = 59999
$ws.Range("N56").ClearContents()
$ws.Range("J82").Value = 2000
$ws.Range("L82").Value = 2000
$ws.Range("H82").Value = 2000
$ws.Range("N82").Value = -2722
$ws.Range("J85").Value = 2000
$ws.Range("L85").Value = 2000
$ws.Range("H85").Value = 2000
$ws.Range("N85").Value = -4496
$ws.Range("J122").Value = 4833.3335
$ws.Range("I122").Value = 3501.3333
$ws.Range("M122").Value = -8053.999899999999
$ws.Range("K122").Value = 10503.9999
$ws.Range("L122").Value = 14500.0005
$ws.Range("H122").Value = 4167.3335
$ws.Range("N122").Value = -19400.0005
$ws.Range("I132").Value = 19502
$ws.Range("M132").Value = -55976
$ws.Range("K132").Value = 58506
$ws.Range("H132").Value = 16334.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J61").Value = 0
$ws.Range("I61").Value = 30486.75
$ws.Range("M61").Value = -30194.75
$ws.Range("K61").Value = 30486.75
$ws.Range("L61").Value = 0
$ws.Range("H61").Value = 30486.75
$ws.Range("N61").ClearContents()
$ws.Range("J70").Value = 45000
$ws.Range("I70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 45000
$ws.Range("H70").Value = 45000
$ws.Range("N70").Value = -45630
$ws.Range("J73").Value = 45000
$ws.Range("I73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 45000
$ws.Range("H73").Value = 45000
$ws.Range("N73").Value = -47184
$ws.Range("I81").Value = 721
$ws.Range("M81").Value = -381
$ws.Range("K81").Value = 1442
$ws.Range("H81").Value = 766.8
$ws.Range("I84").Value = 721
$ws.Range("M84").Value = -1906
$ws.Range("K84").Value = 7210
$ws.Range("H84").Value = 766.8
$ws.Range("J112").Value = 48462.5
$ws.Range("L112").Value = 48462.5
$ws.Range("H112").Value = 48462.5
$ws.Range("N112").Value = -51416.5
$ws.Range("J136").Value = 2123.5
$ws.Range("I136").Value = 2142
$ws.Range("M136").Value = -3876
$ws.Range("K136").Value = 6426
$ws.Range("L136").Value = 6370.5
$ws.Range("H136").Value = 2138.3
$ws.Range("N136").Value = -11470.5
